# Fill in the "0-shot" results row for the three "Prompt_id 1" sub-tables
# (Results @ 1, Results @ 2, Results @ 5) on Sheet1. These rows (77, 82, 87)
# previously only had the row label (A) and blank, pre-styled percent cells
# (E, G); this commit ("results for prompt1 0-shot") populates them with the
# actual measured counts and rates, entered as plain values (not formulas),
# matching how the author pasted the numbers in by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prompt_id 1 - Results @ 1 (row 77)
$ws.Range("B77").Value = 1
$ws.Range("C77").Value = 11
$ws.Range("D77").Value = 93
$ws.Range("E77").Value = 0.0095200000000000007
$ws.Range("F77").Value = 93
$ws.Range("G77").Value = 0.88571

# Prompt_id 1 - Results @ 2 (row 82)
$ws.Range("B82").Value = 6
$ws.Range("C82").Value = 14
$ws.Range("D82").Value = 85
$ws.Range("E82").Value = 0.057140000000000003
$ws.Range("F82").Value = 189
$ws.Range("G82").Value = 0.9

# Prompt_id 1 - Results @ 5 (row 87)
$ws.Range("B87").Value = 10
$ws.Range("C87").Value = 23
$ws.Range("D87").Value = 72
$ws.Range("E87").Value = 0.095000000000000001
$ws.Range("F87").Value = 480
$ws.Range("G87").Value = 0.91429000000000005

# Leave the selection where the author ended up after typing the numbers in.
$ws.Range("G88").Select()
